# Edit: replace the "janellemonae3" tracklist (Sheet1/Sheet3 data range A1:E15)
# with the "cardib1" tracklist (A1:E14), update the two scoped defined names
# and query-table names from janellemonae3 -> cardib1, adjust column widths,
# and update Sheet2's selection.

$wb = $excel.ActiveWorkbook

# New track data (title, composer, performer, time-of-day fraction)
$rows = @(
    @("Get Up 10", "Sean Allen / Belcalis Almanzar / Maurice Jordan / Jermaine Preyan / James SwanQo / Anthony Tucker / Robert Williams", "Cardi B feat. Migos", 0.16041666666666668),
    @("Drip", "Belcalis Almanzar / Kirsnick Ball / Kiari Cephus / Joshua Cross / Quavious Marshall", "Cardi B", 0.18263888888888891),
    @("Bickenhead", "Belcalis Almanzar / James Foye III / Philip Coleman, Jr. / Austin Owens / Jordan Thorpe", "Cardi B", 0.12569444444444444),
    @("Bodak Yellow", "Belcalis Almanzar / Laquan Green / J. White Did It / Dieuson Octave / Klenord Raphael / Jordan Thorpe / Jermaine White", "Cardi B", 0.15486111111111112),
    @("Be Careful", "Belcalis Almanzar / Alan Bergman / Marilyn Bergman / Dennis Coles / Robert Diggs / Adam Feeney / Gary Grice / Marvin Hamlisch / Lamont Hawkins / Anderson Hernandez / Lauryn Hill / Jason Hunter / Russell Jones / Matthew Samuels / Clifford Smith / Jordan Thorpe / U-God / Corey Woods", "Cardi B", 0.14583333333333334),
    @("Best Life", "Belcalis Almanzar / Chancelor Bennett / Allen Ritter / Matthew Samuels", "Cardi B feat. Chacne The Rapper", 0.19722222222222222),
    @("I Like It", "Belcalis Almanzar / Jos? ?lvaro Osorio Balvin / Benito Antonio Martinez Ocasio / Tony Pabon / Manny Rodriguez", "Cardi B feat. Bad Bunny, J. Balvin", 0.17569444444444446),
    @("Ring", "Belcalis Almanzar / Khari Cain / Nija Charles / Kehlani Parrish / Mike Riley", "Cardi B feat. Kehlani", 0.12291666666666667),
    @("Money Bag", "Belcalis Almanzar / Jordan Thorpe / Jermaine White", "Cardi B", 0.15902777777777777),
    @("Bartier Cardi", "Shayaa Bin Abraham-Joseph / Belcalis Almanzar / Samuel Gloade / Darryl McCorkell", "Cardi B feat, 21 Savage", 0.15555555555555556),
    @("She Bad", "Belcalis Almanzar / Keenon Jackson / Leslie Andre Wakefield Jr. / Dijon McFarlane", "Cardi B feat, YG", 0.15972222222222224),
    @("Thru Your Phone", "Belcalis Almanzar / Benjamin Levin / Alexandra Tamposi / Justin Tranter / Andrew Wotman", "Cardi B", 0.13055555555555556),
    @("I Do", "Belcalis Almanzar / Nija Charles / Kevin Gomringer / Tim Gomringer / Shane Lindstrom / Solana Rowe", "Cardi B feat. SZA", 0.1388888888888889)
)

$sheetNames = @("Sheet1", "Sheet3")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = 2 + $i
        $rowData = $rows[$i]
        $ws.Cells.Item($r, 1).Value = ($i + 1)
        $ws.Cells.Item($r, 2).Value = $rowData[0]
        $ws.Cells.Item($r, 3).Value = $rowData[1]
        $ws.Cells.Item($r, 4).Value = $rowData[2]
        $ws.Cells.Item($r, 5).Value = $rowData[3]
    }

    # Row 15 (old 14th track) is no longer used; clear its contents but keep
    # the formatted (empty) E15 cell.
    $ws.Range("A15:E15").ClearContents()

    # Column width tweaks: B 17.77734375 -> 16, D 34.21875 -> 31
    $ws.Columns.Item(2).ColumnWidth = 15.285714285714286
    $ws.Columns.Item(4).ColumnWidth = 30.285714285714285
}

# Update the two sheet-scoped defined names (range shrinks from E15 to E14,
# and the name itself changes from janellemonae3 -> cardib1).
$name1 = $wb.Names.Item(1)
$name1.RefersTo = "=Sheet1!`$A`$1:`$E`$14"
$name2 = $wb.Names.Item(2)
$name2.RefersTo = "=Sheet3!`$A`$1:`$E`$14"
$name1.Name = "cardib1"
$name2.Name = "cardib1"

# Sheet2 selection moves from A3:K18 (active K18) to A3:K3 (active K3).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A3:K3").Select()
